# Generate Report for Handoff
#
# Adds two new handoff records to the Overview, zh-cn and de-de sheets /
# tables:
#   ca41bb4e-c8ba-4dca-acb7-30cca250a57a.md
#   fef6ae5c-543a-4f3a-9059-e029212d393f.md

$wb = $excel.ActiveWorkbook

$dateFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# New file identifiers / metadata
# ---------------------------------------------------------------------------
$file1Id   = "ca41bb4e-c8ba-4dca-acb7-30cca250a57a"
$file1Name = "$file1Id.md"
$file1Path = "e2e\$file1Id.md"
$file1HoDateZh = "2016-08-13 00:49:53"
$file1HoDateDe = "2016-08-13 00:50:03"
$file1XliffZh = "$file1Id.08b73e5db8f3abcb916120ebdbcbe6bf8ea56afa.zh-cn.xlf"
$file1XliffDe = "$file1Id.08b73e5db8f3abcb916120ebdbcbe6bf8ea56afa.de-de.xlf"
$file1Url = "https://github.com/OpenLocalizationTestOrg/oltest/blob/ca41bb4e1f4869cf9033a06fa549f79329ddca41/e2e/$file1Name"

$file2Id   = "fef6ae5c-543a-4f3a-9059-e029212d393f"
$file2Name = "$file2Id.md"
$file2Path = "e2e\$file2Id.md"
$file2HoDateZh = "2016-08-13 00:49:53"
$file2HoDateDe = "2016-08-13 00:50:03"
$file2XliffZh = "$file2Id.714b4a3bc4d10a7acdb38f73e096eaf0cc23f49a.zh-cn.xlf"
$file2XliffDe = "$file2Id.714b4a3bc4d10a7acdb38f73e096eaf0cc23f49a.de-de.xlf"
$file2Url = "https://github.com/OpenLocalizationTestOrg/oltest/blob/fef6ae5c01613719fe7d324877418ead3c2fef6/e2e/$file2Name"

$status = "Ready for handoff"
$extension = ".md"
$naDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Overview sheet
#   A File Name | B Path And Name | C Extension | D Publish URL
#   E zh-cn | F de-de | G Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $file1Name
$wsOverview.Range("B4").Value = $file1Path
$wsOverview.Range("C4").Value = $extension
$wsOverview.Range("D4").Value = ""
$wsOverview.Range("E4").Value = $status
$wsOverview.Range("F4").Value = $status
$wsOverview.Range("G4").Value = $file1HoDateDe
$wsOverview.Range("G4").NumberFormat = $dateFormat
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), $file1Url, "", "", $file1Path) | Out-Null

$wsOverview.Range("A5").Value = $file2Name
$wsOverview.Range("B5").Value = $file2Path
$wsOverview.Range("C5").Value = $extension
$wsOverview.Range("D5").Value = ""
$wsOverview.Range("E5").Value = $status
$wsOverview.Range("F5").Value = $status
$wsOverview.Range("G5").Value = $file2HoDateDe
$wsOverview.Range("G5").NumberFormat = $dateFormat
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), $file2Url, "", "", $file2Path) | Out-Null

$tblOverview = $wsOverview.ListObjects.Item(1)
$tblOverview.Resize($wsOverview.Range("A1:G5"))

# ---------------------------------------------------------------------------
# zh-cn sheet
#   A Source File Name | B File Extension | C Status | D Source Path
#   E Priority | F Content Duplicate | G Latest Handoff File
#   H Latest Handoff Datetime | I Latest Target File | J Latest Handback File
#   K Latest Handback DateTime | L Reference Tokens | M To be localized
#   N Dependency From | O Has metadata | P Error Detail
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = $file1Name
$wsZh.Range("B4").Value = $extension
$wsZh.Range("C4").Value = $status
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'False"
$wsZh.Range("G4").Value = $file1XliffZh
$wsZh.Range("H4").Value = $file1HoDateZh
$wsZh.Range("H4").NumberFormat = $dateFormat
$wsZh.Range("I4").Value = ""
$wsZh.Range("J4").Value = ""
$wsZh.Range("K4").Value = $naDate
$wsZh.Range("K4").NumberFormat = $dateFormat
$wsZh.Range("L4").Value = ""
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = ""
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = ""
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $file1Url, "", "", $file1Name) | Out-Null

$wsZh.Range("A5").Value = $file2Name
$wsZh.Range("B5").Value = $extension
$wsZh.Range("C5").Value = $status
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = "'False"
$wsZh.Range("G5").Value = $file2XliffZh
$wsZh.Range("H5").Value = $file2HoDateZh
$wsZh.Range("H5").NumberFormat = $dateFormat
$wsZh.Range("I5").Value = ""
$wsZh.Range("J5").Value = ""
$wsZh.Range("K5").Value = $naDate
$wsZh.Range("K5").NumberFormat = $dateFormat
$wsZh.Range("L5").Value = ""
$wsZh.Range("M5").Value = "'True"
$wsZh.Range("N5").Value = ""
$wsZh.Range("O5").Value = "'False"
$wsZh.Range("P5").Value = ""
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), $file2Url, "", "", $file2Name) | Out-Null

$tblZh = $wsZh.ListObjects.Item(1)
$tblZh.Resize($wsZh.Range("A1:P5"))

# ---------------------------------------------------------------------------
# de-de sheet (same layout as zh-cn)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = $file1Name
$wsDe.Range("B4").Value = $extension
$wsDe.Range("C4").Value = $status
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'False"
$wsDe.Range("G4").Value = $file1XliffDe
$wsDe.Range("H4").Value = $file1HoDateDe
$wsDe.Range("H4").NumberFormat = $dateFormat
$wsDe.Range("I4").Value = ""
$wsDe.Range("J4").Value = ""
$wsDe.Range("K4").Value = $naDate
$wsDe.Range("K4").NumberFormat = $dateFormat
$wsDe.Range("L4").Value = ""
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = ""
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = ""
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $file1Url, "", "", $file1Name) | Out-Null

$wsDe.Range("A5").Value = $file2Name
$wsDe.Range("B5").Value = $extension
$wsDe.Range("C5").Value = $status
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = "'False"
$wsDe.Range("G5").Value = $file2XliffDe
$wsDe.Range("H5").Value = $file2HoDateDe
$wsDe.Range("H5").NumberFormat = $dateFormat
$wsDe.Range("I5").Value = ""
$wsDe.Range("J5").Value = ""
$wsDe.Range("K5").Value = $naDate
$wsDe.Range("K5").NumberFormat = $dateFormat
$wsDe.Range("L5").Value = ""
$wsDe.Range("M5").Value = "'True"
$wsDe.Range("N5").Value = ""
$wsDe.Range("O5").Value = "'False"
$wsDe.Range("P5").Value = ""
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), $file2Url, "", "", $file2Name) | Out-Null

$tblDe = $wsDe.ListObjects.Item(1)
$tblDe.Resize($wsDe.Range("A1:P5"))
